$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data (rows 2-13), columns A-T:
#   A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
#   E..T = the various expression/specificity metrics.
$data = @(
    @('ECs', 'Efna4', 'Epha2', 'ECs', 2, 0.6666666666666666, 0.6731353333333333, 2.019406, 0.3272865828458516, 0.3272865828458516, 3, 1, 10.004147, 30.012441, 0.4957413074471153, 0.4957413074471153, 6.734144825560668, 60.60730343004601, 0.1622494784899011, 0.1622494784899011),
    @('ECs', 'Efna4', 'Epha2', 'FAPs', 2, 0.6666666666666666, 0.6731353333333333, 2.019406, 0.3272865828458516, 0.3272865828458516, 3, 1, 1.137802, 3.413406, 0.05638216342642133, 0.05638216342642133, 0.7658947285373333, 6.893052556836, 0.01845312560128979, 0.01845312560128979),
    @('ECs', 'Efna4', 'Epha2', 'M2', 2, 0.6666666666666666, 0.6731353333333333, 2.019406, 0.3272865828458516, 0.3272865828458516, 3, 1, 0.5560033333333333, 1.66801, 0.02755195614494878, 0.02755195614494878, 0.3742654891177777, 3.36838940206, 0.009017385577399048, 0.009017385577399048),
    @('ECs', 'Efna4', 'Epha2', 'sCs', 2, 0.6666666666666666, 0.6731353333333333, 2.019406, 0.3272865828458516, 0.3272865828458516, 3, 1, 8.482224, 25.446672, 0.4203245729815145, 0.4203245729815145, 5.709684679648, 51.387162116832, 0.1375665931772617, 0.1375665931772617),
    @('FAPs', 'Efna4', 'Epha2', 'ECs', 3, 1, 0.9964423333333334, 2.989327, 0.4844823769162027, 0.4844823769162026, 3, 1, 10.004147, 30.012441, 0.4957413074471153, 0.4957413074471153, 9.968555579689669, 89.71700021720702, 0.2401779269675244, 0.2401779269675244),
    @('FAPs', 'Efna4', 'Epha2', 'FAPs', 3, 1, 0.9964423333333334, 2.989327, 0.4844823769162027, 0.4844823769162026, 3, 1, 1.137802, 3.413406, 0.05638216342642133, 0.05638216342642133, 1.133754079751333, 10.203786717762, 0.0273161645525104, 0.02731616455251039),
    @('FAPs', 'Efna4', 'Epha2', 'M2', 3, 1, 0.9964423333333334, 2.989327, 0.4844823769162027, 0.4844823769162026, 3, 1, 0.5560033333333333, 1.66801, 0.02755195614494878, 0.02755195614494878, 0.5540252588077778, 4.98622732927, 0.01334843720179576, 0.01334843720179576),
    @('FAPs', 'Efna4', 'Epha2', 'sCs', 3, 1, 0.9964423333333334, 2.989327, 0.4844823769162027, 0.4844823769162026, 3, 1, 8.482224, 25.446672, 0.4203245729815145, 0.4203245729815145, 8.452047074416001, 76.06842366974401, 0.2036398481943721, 0.203639848194372),
    @('sCs', 'Efna4', 'Epha2', 'ECs', 3, 1, 0.3871376666666667, 1.161413, 0.1882310402379457, 0.1882310402379457, 3, 1, 10.004147, 30.012441, 0.4957413074471153, 0.4957413074471153, 3.872982126570334, 34.856839139133, 0.09331390198968979, 0.09331390198968979),
    @('sCs', 'Efna4', 'Epha2', 'FAPs', 3, 1, 0.3871376666666667, 1.161413, 0.1882310402379457, 0.1882310402379457, 3, 1, 1.137802, 3.413406, 0.05638216342642133, 0.05638216342642133, 0.4404860114086667, 3.964374102678, 0.01061287327262114, 0.01061287327262114),
    @('sCs', 'Efna4', 'Epha2', 'M2', 3, 1, 0.3871376666666667, 1.161413, 0.1882310402379457, 0.1882310402379457, 3, 1, 0.5560033333333333, 1.66801, 0.02755195614494878, 0.02755195614494878, 0.2152498331255555, 1.93724849813, 0.00518613336575397, 0.00518613336575397),
    @('sCs', 'Efna4', 'Epha2', 'sCs', 3, 1, 0.3871376666666667, 1.161413, 0.1882310402379457, 0.1882310402379457, 3, 1, 8.482224, 25.446672, 0.4203245729815145, 0.4203245729815145, 3.283788407504, 29.554095667536, 0.07911813160988082, 0.07911813160988082)
)

$r = 2
foreach ($row in $data) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r = $r + 1
}
